# Weekly update: add the latest week's Membrillo price data
# (Fecha = 2022-04-21, Excel serial 44672) for "Terminal La Palmera de
# La Serena" - this pushes the previously existing data rows down by
# two rows (old row 11 -> new row 13, ... old row 37 -> new row 39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the first data row (row 11),
# shifting the rest of the table (old rows 11-37) down to rows 13-39.
# Excel's native row-insert behaviour copies the formatting (e.g. the
# date number format on column D) from the row above, exactly as it
# would via the UI / real COM automation.
$ws.Rows("11:12").Insert()

# Column headers (for reference):
# A Mercado ID | B Mercado | C Región | D Fecha | E Codreg | F Tipo
# G Producto ID | H Producto | I Categoría ID | J Categoría | K Variedad
# L Calidad | M Volumen | N Precio mínimo | O Precio máximo
# P Precio promedio ponderado | Q Unidad de comercialización
# R Origen | S Precio $/Kg | T Kg / unidad

# New row 11: Especial
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Terminal La Palmera de La Serena"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44672
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100104
$ws.Range("H11").Value = "Frutos de pepita"
$ws.Range("I11").Value = 100104003
$ws.Range("J11").Value = "Membrillo"
$ws.Range("K11").Value = "Champion"
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 20
$ws.Range("N11").Value = 320000
$ws.Range("O11").Value = 330000
$ws.Range("P11").Value = 325000
$ws.Range("Q11").Value = "`$/bins (450 kilos)"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 722
$ws.Range("T11").Value = 450

# New row 12: Primera
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = "Terminal La Palmera de La Serena"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44672
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100104
$ws.Range("H12").Value = "Frutos de pepita"
$ws.Range("I12").Value = 100104003
$ws.Range("J12").Value = "Membrillo"
$ws.Range("K12").Value = "Champion"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 10
$ws.Range("N12").Value = 280000
$ws.Range("O12").Value = 290000
$ws.Range("P12").Value = 285000
$ws.Range("Q12").Value = "`$/bins (450 kilos)"
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 633
$ws.Range("T12").Value = 450
